$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (84, 85) following the existing pattern in the tracker.
$ws.Range("A84").Value = "G1"
$ws.Range("B84").Value = "Test1"
$ws.Range("C84").Value = 45902
$ws.Range("C84").NumberFormat = "YYYY-MM-DD"
$ws.Range("D84").Value = 0.6716531388604381
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = -0.01

$ws.Range("A85").Value = "G2"
$ws.Range("B85").Value = "sedrftgyhuioygtfrd"
$ws.Range("C85").Value = 45902
$ws.Range("C85").NumberFormat = "YYYY-MM-DD"
$ws.Range("D85").Value = 0.6716531388604381
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = -0.01
